$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update quantities in row 2 (SERINGAS, LUVAS, TUBOS_DE_COLETA now all 20)
$ws.Range("B2").Value = 20
$ws.Range("D2").Value = 20
$ws.Range("E2").Value = 20

# Select cell J5 and apply an underline font, as the last recorded command
$ws.Range("J5").Select()
$ws.Range("J5").Font.Underline = $true
